$d = $word.ActiveDocument

$d.Content.Find.Execute("60+13=73", $true, $false, $false, $false, $false, $false, 1, $false, "5+31=36", 2) | Out-Null
$d.Content.Find.Execute("7+65=72", $true, $false, $false, $false, $false, $false, 1, $false, "72+7=79", 2) | Out-Null
$d.Content.Find.Execute("16+77=93", $true, $false, $false, $false, $false, $false, 1, $false, "71+13=84", 2) | Out-Null
$d.Content.Find.Execute("96-90=6", $true, $false, $false, $false, $false, $false, 1, $false, "68-24=44", 2) | Out-Null
$d.Content.Find.Execute("35+29=64", $true, $false, $false, $false, $false, $false, 1, $false, "39+12=51", 2) | Out-Null
$d.Content.Find.Execute("2+16=18", $true, $false, $false, $false, $false, $false, 1, $false, "81-18=63", 2) | Out-Null
$d.Content.Find.Execute("27+15=42", $true, $false, $false, $false, $false, $false, 1, $false, "51-4=47", 2) | Out-Null
$d.Content.Find.Execute("45-17=28", $true, $false, $false, $false, $false, $false, 1, $false, "0+34=34", 2) | Out-Null
$d.Content.Find.Execute("44-38=6", $true, $false, $false, $false, $false, $false, 1, $false, "8+48=56", 2) | Out-Null
$d.Content.Find.Execute("61-56=5", $true, $false, $false, $false, $false, $false, 1, $false, "23+72=95", 2) | Out-Null
$d.Content.Find.Execute("5+39=44", $true, $false, $false, $false, $false, $false, 1, $false, "4+92=96", 2) | Out-Null
$d.Content.Find.Execute("22+57=79", $true, $false, $false, $false, $false, $false, 1, $false, "5+83=88", 2) | Out-Null
$d.Content.Find.Execute("37+55=92", $true, $false, $false, $false, $false, $false, 1, $false, "4+14=18", 2) | Out-Null
$d.Content.Find.Execute("65+0=65", $true, $false, $false, $false, $false, $false, 1, $false, "26-14=12", 2) | Out-Null
$d.Content.Find.Execute("75-72=3", $true, $false, $false, $false, $false, $false, 1, $false, "77-44=33", 2) | Out-Null
$d.Content.Find.Execute("10+51=61", $true, $false, $false, $false, $false, $false, 1, $false, "13+0=13", 2) | Out-Null
$d.Content.Find.Execute("14+38=52", $true, $false, $false, $false, $false, $false, 1, $false, "28+38=66", 2) | Out-Null
$d.Content.Find.Execute("41+6=47", $true, $false, $false, $false, $false, $false, 1, $false, "16+48=64", 2) | Out-Null
$d.Content.Find.Execute("82+7=89", $true, $false, $false, $false, $false, $false, 1, $false, "18-15=3", 2) | Out-Null
$d.Content.Find.Execute("79-0=79", $true, $false, $false, $false, $false, $false, 1, $false, "26+62=88", 2) | Out-Null
$d.Content.Find.Execute("45+12=57", $true, $false, $false, $false, $false, $false, 1, $false, "66-29=37", 2) | Out-Null
$d.Content.Find.Execute("54+35=89", $true, $false, $false, $false, $false, $false, 1, $false, "1+23=24", 2) | Out-Null
$d.Content.Find.Execute("32+2=34", $true, $false, $false, $false, $false, $false, 1, $false, "21+63=84", 2) | Out-Null
$d.Content.Find.Execute("77-53=24", $true, $false, $false, $false, $false, $false, 1, $false, "28-21=7", 2) | Out-Null
$d.Content.Find.Execute("58+27=85", $true, $false, $false, $false, $false, $false, 1, $false, "88-11=77", 2) | Out-Null
$d.Content.Find.Execute("21+4=25", $true, $false, $false, $false, $false, $false, 1, $false, "66+26=92", 2) | Out-Null
$d.Content.Find.Execute("49+10=59", $true, $false, $false, $false, $false, $false, 1, $false, "37-5=32", 2) | Out-Null
$d.Content.Find.Execute("17+16=33", $true, $false, $false, $false, $false, $false, 1, $false, "54+11=65", 2) | Out-Null
$d.Content.Find.Execute("67-58=9", $true, $false, $false, $false, $false, $false, 1, $false, "16+61=77", 2) | Out-Null
$d.Content.Find.Execute("46+30=76", $true, $false, $false, $false, $false, $false, 1, $false, "91-12=79", 2) | Out-Null
$d.Content.Find.Execute("70-33=37", $true, $false, $false, $false, $false, $false, 1, $false, "60-40=20", 2) | Out-Null
$d.Content.Find.Execute("3+76=79", $true, $false, $false, $false, $false, $false, 1, $false, "59+12=71", 2) | Out-Null
$d.Content.Find.Execute("64+32=96", $true, $false, $false, $false, $false, $false, 1, $false, "53-43=10", 2) | Out-Null
$d.Content.Find.Execute("72-20=52", $true, $false, $false, $false, $false, $false, 1, $false, "81-9=72", 2) | Out-Null
$d.Content.Find.Execute("92-50=42", $true, $false, $false, $false, $false, $false, 1, $false, "91-73=18", 2) | Out-Null
$d.Content.Find.Execute("2+79=81", $true, $false, $false, $false, $false, $false, 1, $false, "29+29=58", 2) | Out-Null
$d.Content.Find.Execute("5+2=7", $true, $false, $false, $false, $false, $false, 1, $false, "55+30=85", 2) | Out-Null
$d.Content.Find.Execute("22-6=16", $true, $false, $false, $false, $false, $false, 1, $false, "35+10=45", 2) | Out-Null
$d.Content.Find.Execute("47-22=25", $true, $false, $false, $false, $false, $false, 1, $false, "80-77=3", 2) | Out-Null
$d.Content.Find.Execute("88-86=2", $true, $false, $false, $false, $false, $false, 1, $false, "91-78=13", 2) | Out-Null
$d.Content.Find.Execute("22+68=90", $true, $false, $false, $false, $false, $false, 1, $false, "29+26=55", 2) | Out-Null
$d.Content.Find.Execute("13+26=39", $true, $false, $false, $false, $false, $false, 1, $false, "34+11=45", 2) | Out-Null
$d.Content.Find.Execute("95-60=35", $true, $false, $false, $false, $false, $false, 1, $false, "2+96=98", 2) | Out-Null
$d.Content.Find.Execute("7+70=77", $true, $false, $false, $false, $false, $false, 1, $false, "1+14=15", 2) | Out-Null
$d.Content.Find.Execute("93-88=5", $true, $false, $false, $false, $false, $false, 1, $false, "87+6=93", 2) | Out-Null
$d.Content.Find.Execute("4+77=81", $true, $false, $false, $false, $false, $false, 1, $false, "8-5=3", 2) | Out-Null
$d.Content.Find.Execute("94-36=58", $true, $false, $false, $false, $false, $false, 1, $false, "85+13=98", 2) | Out-Null
$d.Content.Find.Execute("52+38=90", $true, $false, $false, $false, $false, $false, 1, $false, "93-7=86", 2) | Out-Null
$d.Content.Find.Execute("75-29=46", $true, $false, $false, $false, $false, $false, 1, $false, "52-3=49", 2) | Out-Null
$d.Content.Find.Execute("93-28=65", $true, $false, $false, $false, $false, $false, 1, $false, "21+70=91", 2) | Out-Null
$d.Content.Find.Execute("29+5=34", $true, $false, $false, $false, $false, $false, 1, $false, "7+50=57", 2) | Out-Null
$d.Content.Find.Execute("18+70=88", $true, $false, $false, $false, $false, $false, 1, $false, "16+23=39", 2) | Out-Null
$d.Content.Find.Execute("52-48=4", $true, $false, $false, $false, $false, $false, 1, $false, "47+20=67", 2) | Out-Null
$d.Content.Find.Execute("12+19=31", $true, $false, $false, $false, $false, $false, 1, $false, "31-15=16", 2) | Out-Null
$d.Content.Find.Execute("9+35=44", $true, $false, $false, $false, $false, $false, 1, $false, "12+34=46", 2) | Out-Null
$d.Content.Find.Execute("44-3=41", $true, $false, $false, $false, $false, $false, 1, $false, "34-6=28", 2) | Out-Null
$d.Content.Find.Execute("95-7=88", $true, $false, $false, $false, $false, $false, 1, $false, "89-26=63", 2) | Out-Null
$d.Content.Find.Execute("65-61=4", $true, $false, $false, $false, $false, $false, 1, $false, "96-3=93", 2) | Out-Null
$d.Content.Find.Execute("85+2=87", $true, $false, $false, $false, $false, $false, 1, $false, "44-34=10", 2) | Out-Null
$d.Content.Find.Execute("98-98=0", $true, $false, $false, $false, $false, $false, 1, $false, "30+7=37", 2) | Out-Null
$d.Content.Find.Execute("89-41=48", $true, $false, $false, $false, $false, $false, 1, $false, "80-76=4", 2) | Out-Null
$d.Content.Find.Execute("30+2=32", $true, $false, $false, $false, $false, $false, 1, $false, "61-58=3", 2) | Out-Null
$d.Content.Find.Execute("19+37=56", $true, $false, $false, $false, $false, $false, 1, $false, "52-31=21", 2) | Out-Null
$d.Content.Find.Execute("99-44=55", $true, $false, $false, $false, $false, $false, 1, $false, "16+71=87", 2) | Out-Null
$d.Content.Find.Execute("4+94=98", $true, $false, $false, $false, $false, $false, 1, $false, "10+68=78", 2) | Out-Null
$d.Content.Find.Execute("3+58=61", $true, $false, $false, $false, $false, $false, 1, $false, "32+21=53", 2) | Out-Null
$d.Content.Find.Execute("10+55=65", $true, $false, $false, $false, $false, $false, 1, $false, "17+17=34", 2) | Out-Null
$d.Content.Find.Execute("77-64=13", $true, $false, $false, $false, $false, $false, 1, $false, "76-42=34", 2) | Out-Null
$d.Content.Find.Execute("94-20=74", $true, $false, $false, $false, $false, $false, 1, $false, "64-12=52", 2) | Out-Null
$d.Content.Find.Execute("45+30=75", $true, $false, $false, $false, $false, $false, 1, $false, "17+10=27", 2) | Out-Null
$d.Content.Find.Execute("16+7=23", $true, $false, $false, $false, $false, $false, 1, $false, "94-19=75", 2) | Out-Null
$d.Content.Find.Execute("99-38=61", $true, $false, $false, $false, $false, $false, 1, $false, "48-23=25", 2) | Out-Null
$d.Content.Find.Execute("96-41=55", $true, $false, $false, $false, $false, $false, 1, $false, "23+9=32", 2) | Out-Null
$d.Content.Find.Execute("40+13=53", $true, $false, $false, $false, $false, $false, 1, $false, "61-34=27", 2) | Out-Null
$d.Content.Find.Execute("55+13=68", $true, $false, $false, $false, $false, $false, 1, $false, "39-3=36", 2) | Out-Null
$d.Content.Find.Execute("32+51=83", $true, $false, $false, $false, $false, $false, 1, $false, "74-48=26", 2) | Out-Null
$d.Content.Find.Execute("73-59=14", $true, $false, $false, $false, $false, $false, 1, $false, "87+11=98", 2) | Out-Null
$d.Content.Find.Execute("88-30=58", $true, $false, $false, $false, $false, $false, 1, $false, "94-12=82", 2) | Out-Null
$d.Content.Find.Execute("26+13=39", $true, $false, $false, $false, $false, $false, 1, $false, "29+59=88", 2) | Out-Null
$d.Content.Find.Execute("51-18=33", $true, $false, $false, $false, $false, $false, 1, $false, "7+43=50", 2) | Out-Null
$d.Content.Find.Execute("36+40=76", $true, $false, $false, $false, $false, $false, 1, $false, "36+34=70", 2) | Out-Null
$d.Content.Find.Execute("11+1=12", $true, $false, $false, $false, $false, $false, 1, $false, "35+49=84", 2) | Out-Null
$d.Content.Find.Execute("36+62=98", $true, $false, $false, $false, $false, $false, 1, $false, "23+56=79", 2) | Out-Null
$d.Content.Find.Execute("23-18=5", $true, $false, $false, $false, $false, $false, 1, $false, "88-3=85", 2) | Out-Null
$d.Content.Find.Execute("98-93=5", $true, $false, $false, $false, $false, $false, 1, $false, "76+2=78", 2) | Out-Null
$d.Content.Find.Execute("14+19=33", $true, $false, $false, $false, $false, $false, 1, $false, "48+51=99", 2) | Out-Null
$d.Content.Find.Execute("11+8=19", $true, $false, $false, $false, $false, $false, 1, $false, "63-1=62", 2) | Out-Null
$d.Content.Find.Execute("10+3=13", $true, $false, $false, $false, $false, $false, 1, $false, "61-23=38", 2) | Out-Null
$d.Content.Find.Execute("96-52=44", $true, $false, $false, $false, $false, $false, 1, $false, "73-31=42", 2) | Out-Null
$d.Content.Find.Execute("20-13=7", $true, $false, $false, $false, $false, $false, 1, $false, "29+45=74", 2) | Out-Null
$d.Content.Find.Execute("55-0=55", $true, $false, $false, $false, $false, $false, 1, $false, "91-42=49", 2) | Out-Null
$d.Content.Find.Execute("22-15=7", $true, $false, $false, $false, $false, $false, 1, $false, "75-22=53", 2) | Out-Null
$d.Content.Find.Execute("51-29=22", $true, $false, $false, $false, $false, $false, 1, $false, "13+14=27", 2) | Out-Null
$d.Content.Find.Execute("31+30=61", $true, $false, $false, $false, $false, $false, 1, $false, "56+20=76", 2) | Out-Null
$d.Content.Find.Execute("63-4=59", $true, $false, $false, $false, $false, $false, 1, $false, "27+3=30", 2) | Out-Null
$d.Content.Find.Execute("23+52=75", $true, $false, $false, $false, $false, $false, 1, $false, "72-5=67", 2) | Out-Null
$d.Content.Find.Execute("28-7=21", $true, $false, $false, $false, $false, $false, 1, $false, "66-40=26", 2) | Out-Null
$d.Content.Find.Execute("15-14=1", $true, $false, $false, $false, $false, $false, 1, $false, "52-6=46", 2) | Out-Null
$d.Content.Find.Execute("42-39=3", $true, $false, $false, $false, $false, $false, 1, $false, "89-5=84", 2) | Out-Null
$d.Content.Find.Execute("92-85=7", $true, $false, $false, $false, $false, $false, 1, $false, "49-16=33", 2) | Out-Null
